# Logs and extent reports fixed.
#
# Rewrites the "common info" sheet (sheet1) so it becomes a generic Key/Value
# table (adding a new "home_page_title" entry with the Takealot home page
# title, keeping "alert_wait_time", keeping "signup_page_title" but blanking
# its value) and drops the old "global_wait_time" / "login_page_title" /
# "product_category" rows. The previously active "search data" tab is
# deselected in favour of "common info".

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "common info"

# --- Rewrite sheet1 content -------------------------------------------------
# Order of writes matters for where new shared strings land in the table, so
# the new unique string "Takealot.com: ..." is entered before the "Key" /
# "Value" header labels and the "home_page_title" label.
$ws1.Range("B2").Value = "Takealot.com: Online Shopping | SA's leading online store"
$ws1.Range("A1").Value = "Key"
$ws1.Range("B1").Value = "Value"
$ws1.Range("A2").Value = "home_page_title"
$ws1.Range("A1:B1").Font.Bold = $true

$ws1.Range("A3").Value = "alert_wait_time"
$ws1.Range("B3").Value = 10

$ws1.Range("A4").Value = "signup_page_title"
$ws1.Range("B4").ClearContents()

# Drop the old trailing row (former "product_category" / "phones" row) so the
# sheet shrinks from 5 rows to 4.
$ws1.Rows.Item(5).Delete()

# Make "common info" the active sheet/selection (it replaces "search data" as
# the tab shown when the workbook is opened).
[void]$ws1.Range("A2").Select()

Write-Host "done"
